$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-22 Sunday" "2024-12-23 Monday"

Replace-Text "693×6=4158" "943×3=2829"
Replace-Text "714×2=1428" "564×4=2256"
Replace-Text "832×4=3328" "897×5=4485"
Replace-Text "448×9=4032" "227×3=681"
Replace-Text "313×7=2191" "394×4=1576"

Replace-Text "274×2=548" "268×5=1340"
Replace-Text "620×8=4960" "823×4=3292"
Replace-Text "781×2=1562" "976×7=6832"
Replace-Text "529×8=4232" "869×4=3476"
Replace-Text "955×9=8595" "551×9=4959"

Replace-Text "461×6=2766" "408×7=2856"
Replace-Text "673×4=2692" "218×8=1744"
Replace-Text "302×6=1812" "293×3=879"
Replace-Text "791×7=5537" "456×4=1824"
Replace-Text "843×7=5901" "693×9=6237"

Replace-Text "674×8=5392" "911×4=3644"
Replace-Text "211×9=1899" "375×5=1875"
Replace-Text "526×8=4208" "573×6=3438"
Replace-Text "305×2=610" "973×6=5838"
Replace-Text "439×3=1317" "604×7=4228"

Replace-Text "117×4=468" "899×4=3596"
Replace-Text "771×2=1542" "831×5=4155"
Replace-Text "468×2=936" "135×8=1080"
Replace-Text "525×5=2625" "252×6=1512"
Replace-Text "431×7=3017" "178×6=1068"
